# Generated with Celerio v3.0.101 - http://www.springfuse.com/
#
# This edit reworks the "book" entity fields used by the generated Excel
# export templates:
#   - book_accountId / book.accountId  ->  book_owner / printer.print(book.owner)
#   - book_title      / book.title      ->  book_bookTitle / book.bookTitle
# and updates the "Search" sheet so that the single "title" search box is
# split into two separate search rows: one for "owner" and one for
# "bookTitle" (pushing the existing numberOfPages range-search row down
# by one row).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List")
$ws2 = $wb.Worksheets.Item("Search")

# ---------------------------------------------------------------------
# Sheet "List": rename the accountId/title columns to owner/bookTitle
# ---------------------------------------------------------------------
$ws1.Range("B1").Value = '${msg.getProperty(''book_owner'')}'
$ws1.Range("B2").Value = '${printer.print(book.owner)}'

$ws1.Range("C1").Value = '${msg.getProperty(''book_bookTitle'')}'
$ws1.Range("C2").Value = '${book.bookTitle}'

# ---------------------------------------------------------------------
# Sheet "Search": split the old "title" search row (row 5) into two rows
# (owner, bookTitle) and shift the numberOfPages range search down to
# row 7.
# ---------------------------------------------------------------------

# Insert a blank row at position 6; this pushes the existing row 6
# (numberOfPages range search) down to row 7, while row 5 is left alone.
$ws2.Rows.Item(6).Insert()

# Row 5: search by owner
$ws2.Range("A5").Value = '${msg.getProperty(''book_owner'')}'
$ws2.Range("B5").Value = '${owner}'

# Row 6 (new): search by bookTitle
$ws2.Range("A6").Value = '${msg.getProperty(''book_bookTitle'')}'
$ws2.Range("B6").Value = '${bookTitle}'

# Row 7 keeps the original numberOfPages range-search content untouched
# (it was shifted down automatically by the row insertion above).
